$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text values are entered in the same order they first appear in the
# workbook's shared-strings table so the resulting index order matches.
$ws.Range("A6").Value = "Server crashes when 'Messages.txt' doesn’t exist and the show messages button is pressed"
$ws.Range("F6").Value = "Yes"
$ws.Range("A10").Value = "Producer crashes when server is offline, need message box to display no message available"
$ws.Range("A7").Value = "Server does not notify user when there are no messages to display"
$ws.Range("A9").Value = "Conumer crashes when requesting a message priority that isnt in the file"
$ws.Range("A8").Value = "Consumer crashes when server is offline, does not notify user"

# Remaining text cells reuse existing shared strings.
$ws.Range("F7").Value = "Yes"
$ws.Range("F8").Value = "Yes"
$ws.Range("F9").Value = "Yes"
$ws.Range("F10").Value = "Yes"

# Numeric values for each new row.
$ws.Range("B6").Value = 5
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0

$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 0

$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 20
$ws.Range("D8").Value = 22
$ws.Range("E8").Value = 0

$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 0

$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 0

$ws.Range("A12").Select()
